$wb = $excel.ActiveWorkbook

# --- Sheet "AR" ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.001117466870434384
$ws.Range("B3").Value = 0.8120860912038282
$ws.Range("B4").Value = 0.102367977347484
$ws.Range("B5").Value = "[1.0, 0.08678975284060744, -0.012840413166194881, -0.07086202935972366, -0.06347007268501068, 0.08296661964641912, -0.0005683124618095891, 0.0545886310286447, 0.03822843924441717, -0.002455064398582463, -0.03561727930752916, -0.01962329612276003, 0.0015551560319931037, 0.06893194384351661, 0.048917368702744074, -0.02678988668141005, -0.02583281521432406, -0.060916076908125165, -0.0368111626394735, -0.02742827854580032]"

# --- Sheet "SETAR" ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.1606068165559258
$ws.Range("B5").Value = 0.6449197196293537
$ws.Range("B6").Value = 0.05683039258314826
$ws.Range("B7").Value = 0.1533935684166628
$ws.Range("B8").Value = 0.6420157588091135
$ws.Range("B9").Value = 0.06860474090610616
$ws.Range("B10").Value = "[1.0, 0.1293755161391712, 0.02645190378559267, -0.011840806256087157, 0.0387866987920049, 0.04207185327433184, 0.0020672235467262727, 0.03799479214923581, 0.008243529135127416, 0.021038990771246702, 0.007484919837517472, -0.015617544388783095, 0.01586141881538815, 0.007963271297386103, 0.0313863721735704, -0.00960337238172656, 0.017302251053871415, -0.028481267567980553, -0.015363432219725354, 0.016200882414218024]"

# --- Sheet "GARCH" ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.002515095599055716
$ws.Range("B3").Value = 0.1018809108678828
$ws.Range("B4").Value = 0.09340073753386273
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.011181391481007911, -0.02352879836637044, -0.08179855948869703, -0.08789542861518558, 0.14071694453960007, 0.021251354655626352, 0.06210508617688904, 0.03462536381979392, -0.023379877313978954, -0.04560155785199741, -0.011497979539804314, 0.00522156820951697, 0.06481729508239266, 0.03894650335553485, -0.0174845378625168, -0.022159819862650365, -0.05175154368543001, -0.038314481793475685, -0.017078761013970555]"

# --- Sheet "TARCH" ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.003223091898555463
$ws.Range("B3").Value = 0.1021758761861158
$ws.Range("B4").Value = 0.1259422174992773
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[0.9999999999999999, 0.011707460101076797, -0.021543155402464077, -0.08124715110605364, -0.08858183052985101, 0.14206606194534827, 0.02432623513627062, 0.05995236962186577, 0.035476217832254534, -0.023589673144808775, -0.04675953136935498, -0.012973206907263923, 0.005195078938895303, 0.0641415500304011, 0.039131509262643714, -0.017944844370046793, -0.020822242105189473, -0.051676254374580954, -0.03836111705883747, -0.017698837057731672]"
$ws.Range("B7").Value = -0.07068099061684387

# --- Sheet "AR_TARCH" ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.0004020974283226129
$ws.Range("B3").Value = 0.09248251131941676
$ws.Range("B4").Value = 0.1259619156021742
$ws.Range("B5").Value = 0.000000000000004475887380093997
$ws.Range("B6").Value = "[1.0, -0.002656763999799318, -0.022497676813174286, -0.07505307542290127, -0.07593642639152301, 0.09108064225172934, -0.003867628583389821, 0.05164712760270787, 0.03217670840622298, -0.007201290129249432, -0.037399915345544765, -0.020838474830338093, -0.005963322968838337, 0.06918348602290686, 0.051656724981900645, -0.026301260207270156, -0.014540574992769388, -0.06104294309570778, -0.03630222563376779, -0.024131294748313638]"
$ws.Range("B7").Value = -0.06526756145371936
$ws.Range("B9").Value = 0.8102371891736481
